$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2935.8096
$ws.Range("I28").Value = 1050.2
$ws.Range("J28").Value = 4650
$ws.Range("K28").Value = 1050.2
$ws.Range("L28").Value = 4650
$ws.Range("M28").Value = -565.2
$ws.Range("N28").Value = -5620

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8851.357
$ws.Range("I51").Value = 50300
$ws.Range("J51").Value = 5663
$ws.Range("K51").Value = 50300
$ws.Range("L51").Value = 5663
$ws.Range("M51").Value = -49816
$ws.Range("N51").Value = -6631

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3671.718
$ws.Range("J112").Value = 3901.5
$ws.Range("L112").Value = 11704.5
$ws.Range("N112").Value = -13920.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 118768.8
$ws.Range("J134").Value = 118768.8
$ws.Range("L134").Value = 118768.8
$ws.Range("N134").Value = -128908.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1339.4286
$ws.Range("I2").Value = 1211.4286
$ws.Range("J2").Value = 1595.4286
$ws.Range("K2").Value = 1211.4286
$ws.Range("L2").Value = 1595.4286
$ws.Range("M2").Value = -1098.4286
$ws.Range("N2").Value = -1821.4286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9265.764999999999
$ws.Range("I32").Value = 9042.0625
$ws.Range("J32").Value = 12845
$ws.Range("K32").Value = 9042.0625
$ws.Range("L32").Value = 12845
$ws.Range("M32").Value = -8755.0625
$ws.Range("N32").Value = -13419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15782984
$ws.Range("I61").Value = 16671948
$ws.Range("J61").Value = 3337499.8
$ws.Range("K61").Value = 16671948
$ws.Range("L61").Value = 3337499.8
$ws.Range("M61").Value = -16671736
$ws.Range("N61").Value = -3337923.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1092
$ws.Range("I97").Value = 733.95
$ws.Range("J97").Value = 2524.2
$ws.Range("K97").Value = 733.95
$ws.Range("L97").Value = 2524.2
$ws.Range("M97").Value = -237.95
$ws.Range("N97").Value = -3516.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1339.4286
$ws.Range("I116").Value = 1211.4286
$ws.Range("J116").Value = 1595.4286
$ws.Range("K116").Value = 1211.4286
$ws.Range("L116").Value = 1595.4286
$ws.Range("M116").Value = 1082.5714
$ws.Range("N116").Value = -6183.4286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3229007.5
$ws.Range("I132").Value = 3166.5862
$ws.Range("K132").Value = 9499.758600000001
$ws.Range("M132").Value = -6969.758600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 15782984
$ws.Range("I136").Value = 16671948
$ws.Range("J136").Value = 3337499.8
$ws.Range("K136").Value = 50015844
$ws.Range("L136").Value = 10012499.4
$ws.Range("M136").Value = -50013294
$ws.Range("N136").Value = -10017599.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1339.4286
$ws.Range("I3").Value = 1211.4286
$ws.Range("J3").Value = 1595.4286
$ws.Range("K3").Value = 1211.4286
$ws.Range("L3").Value = 1595.4286
$ws.Range("M3").Value = -1097.4286
$ws.Range("N3").Value = -1823.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1389
$ws.Range("I86").Value = 1020.82355
$ws.Range("J86").Value = 2640.8
$ws.Range("K86").Value = 1020.82355
$ws.Range("L86").Value = 2640.8
$ws.Range("M86").Value = 102.17645
$ws.Range("N86").Value = -4886.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1389
$ws.Range("I89").Value = 1020.82355
$ws.Range("J89").Value = 2640.8
$ws.Range("K89").Value = 5104.117749999999
$ws.Range("L89").Value = 13204
$ws.Range("M89").Value = 511.8822500000006
$ws.Range("N89").Value = -24436

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1756.5333
$ws.Range("I99").Value = 716.6667
$ws.Range("J99").Value = 3316.3333
$ws.Range("K99").Value = 716.6667
$ws.Range("L99").Value = 3316.3333
$ws.Range("M99").Value = 781.3333
$ws.Range("N99").Value = -6312.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3158.7036
$ws.Range("I107").Value = 3540.9
$ws.Range("J107").Value = 2066.7144
$ws.Range("K107").Value = 3540.9
$ws.Range("L107").Value = 2066.7144
$ws.Range("M107").Value = -1620.9
$ws.Range("N107").Value = -5906.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23258604
$ws.Range("I31").Value = 27029326
$ws.Range("J31").Value = 5824.1665
$ws.Range("K31").Value = 27029326
$ws.Range("L31").Value = 5824.1665
$ws.Range("M31").Value = -27029031
$ws.Range("N31").Value = -6414.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23258604
$ws.Range("I34").Value = 27029326
$ws.Range("J34").Value = 5824.1665
$ws.Range("K34").Value = 27029326
$ws.Range("L34").Value = 5824.1665
$ws.Range("M34").Value = -27029124
$ws.Range("N34").Value = -6228.1665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1637.561
$ws.Range("I134").Value = 1506.0857
$ws.Range("J134").Value = 2404.5
$ws.Range("K134").Value = 4518.257100000001
$ws.Range("L134").Value = 7213.5
$ws.Range("M134").Value = -1983.257100000001
$ws.Range("N134").Value = -12283.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 134
$ws.Range("I13").Value = 134
$ws.Range("K13").Value = 402
$ws.Range("M13").Value = -234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 16359.866
$ws.Range("I56").Value = 16359.866
$ws.Range("K56").Value = 16359.866
$ws.Range("M56").Value = -15829.866

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 7343.1665
$ws.Range("I19").Value = 8473.25
$ws.Range("J19").Value = 5083
$ws.Range("K19").Value = 8473.25
$ws.Range("L19").Value = 5083
$ws.Range("M19").Value = -8185.25
$ws.Range("N19").Value = -5659

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4660380
$ws.Range("I132").Value = 4248.5454
$ws.Range("J132").Value = 23866924
$ws.Range("K132").Value = 12745.6362
$ws.Range("L132").Value = 71600772
$ws.Range("M132").Value = -10215.6362
$ws.Range("N132").Value = -71605832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 22006584
$ws.Range("I22").Value = 26406500
$ws.Range("K22").Value = 26406500
$ws.Range("M22").Value = -26406205

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 22006584
$ws.Range("I27").Value = 26406500
$ws.Range("K27").Value = 26406500
$ws.Range("M27").Value = -26406393

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 61000
$ws.Range("J135").Value = 61000
$ws.Range("L135").Value = 61000
$ws.Range("N135").Value = -71140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 89617.71000000001
$ws.Range("J135").Value = 89617.71000000001
$ws.Range("L135").Value = 89617.71000000001
$ws.Range("N135").Value = -99757.71000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 267120.5
$ws.Range("I136").Value = 4179.8066
$ws.Range("K136").Value = 12539.4198
$ws.Range("M136").Value = -9989.4198
